$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 135, shifting existing rows 135-143 down to 136-144
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new data record
$ws.Cells.Item(135, 1).Value = 3
$ws.Cells.Item(135, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(135, 3).Value = "Coquimbo"
$ws.Cells.Item(135, 4).Value = 44610
$ws.Cells.Item(135, 4).NumberFormat = $ws.Cells.Item(136, 4).NumberFormat
$ws.Cells.Item(135, 5).Value = 5
$ws.Cells.Item(135, 6).Value = 100112052
$ws.Cells.Item(135, 7).Value = "Albahaca"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 85
$ws.Cells.Item(135, 11).Value = 4500
$ws.Cells.Item(135, 12).Value = 5000
$ws.Cells.Item(135, 13).Value = 4735
$ws.Cells.Item(135, 14).Value = "$/docena de matas"
$ws.Cells.Item(135, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(135, 16).Value = 789
$ws.Cells.Item(135, 17).Value = 6
$ws.Cells.Item(135, 18).Value = "Hortaliza"
